$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(3, 8).Value = 13428.286
$ws_ALC.Cells.Item(3, 10).Value = 13428.286
$ws_ALC.Cells.Item(3, 12).Value = 13428.286
$ws_ALC.Cells.Item(3, 14).Value = -13656.286
$ws_ALC.Cells.Item(4, 8).Value = 3665.8333
$ws_ALC.Cells.Item(4, 9).Value = 1726.3636
$ws_ALC.Cells.Item(4, 11).Value = 1726.3636
$ws_ALC.Cells.Item(4, 13).Value = -1612.3636
$ws_ALC.Cells.Item(19, 8).Value = 623.75
$ws_ALC.Cells.Item(19, 10).Value = 641.8570999999999
$ws_ALC.Cells.Item(19, 12).Value = 641.8570999999999
$ws_ALC.Cells.Item(19, 14).Value = -991.8570999999999
$ws_ALC.Cells.Item(32, 8).Value = 26318124
$ws_ALC.Cells.Item(32, 10).Value = 31252090
$ws_ALC.Cells.Item(32, 12).Value = 31252090
$ws_ALC.Cells.Item(32, 14).Value = -31252742
$ws_ALC.Cells.Item(33, 9).Value = 4839188.5
$ws_ALC.Cells.Item(33, 11).Value = 4839188.5
$ws_ALC.Cells.Item(33, 13).Value = -4838959.5
$ws_ALC.Cells.Item(41, 8).Value = 1738.8948
$ws_ALC.Cells.Item(41, 9).Value = 1190.7142
$ws_ALC.Cells.Item(41, 11).Value = 1190.7142
$ws_ALC.Cells.Item(41, 13).Value = -750.7141999999999
$ws_ALC.Cells.Item(43, 8).Value = 54172916
$ws_ALC.Cells.Item(43, 9).Value = 216666670
$ws_ALC.Cells.Item(43, 11).Value = 216666670
$ws_ALC.Cells.Item(43, 13).Value = -216666601
$ws_ALC.Cells.Item(48, 8).Value = 0
$ws_ALC.Cells.Item(48, 10).Value = 0
$ws_ALC.Cells.Item(48, 12).Value = 0
$ws_ALC.Cells.Item(48, 14).ClearContents()
$ws_ALC.Cells.Item(56, 8).Value = 0
$ws_ALC.Cells.Item(56, 10).Value = 0
$ws_ALC.Cells.Item(56, 12).Value = 0
$ws_ALC.Cells.Item(56, 14).ClearContents()
$ws_ALC.Cells.Item(68, 8).Value = 0
$ws_ALC.Cells.Item(68, 10).Value = 0
$ws_ALC.Cells.Item(68, 12).Value = 0
$ws_ALC.Cells.Item(68, 14).ClearContents()
$ws_ALC.Cells.Item(69, 8).Value = 14874.75
$ws_ALC.Cells.Item(69, 9).Value = 14624.75
$ws_ALC.Cells.Item(69, 10).Value = 14999.75
$ws_ALC.Cells.Item(69, 11).Value = 43874.25
$ws_ALC.Cells.Item(69, 12).Value = 44999.25
$ws_ALC.Cells.Item(69, 13).Value = -43000.25
$ws_ALC.Cells.Item(69, 14).Value = -46747.25
$ws_ALC.Cells.Item(70, 8).Value = 1799
$ws_ALC.Cells.Item(70, 10).Value = 0
$ws_ALC.Cells.Item(70, 12).Value = 0
$ws_ALC.Cells.Item(70, 14).ClearContents()
$ws_ALC.Cells.Item(71, 8).Value = 0
$ws_ALC.Cells.Item(71, 10).Value = 0
$ws_ALC.Cells.Item(71, 12).Value = 0
$ws_ALC.Cells.Item(71, 14).ClearContents()
$ws_ALC.Cells.Item(72, 8).Value = 14874.75
$ws_ALC.Cells.Item(72, 9).Value = 14624.75
$ws_ALC.Cells.Item(72, 10).Value = 14999.75
$ws_ALC.Cells.Item(72, 11).Value = 131622.75
$ws_ALC.Cells.Item(72, 12).Value = 134997.75
$ws_ALC.Cells.Item(72, 13).Value = -127254.75
$ws_ALC.Cells.Item(72, 14).Value = -143733.75
$ws_ALC.Cells.Item(73, 8).Value = 1799
$ws_ALC.Cells.Item(73, 10).Value = 0
$ws_ALC.Cells.Item(73, 12).Value = 0
$ws_ALC.Cells.Item(73, 14).ClearContents()
$ws_ALC.Cells.Item(94, 8).Value = 783.2222
$ws_ALC.Cells.Item(94, 9).Value = 783.2222
$ws_ALC.Cells.Item(94, 11).Value = 783.2222
$ws_ALC.Cells.Item(94, 13).Value = -332.2222
$ws_ALC.Cells.Item(99, 8).Value = 397.4
$ws_ALC.Cells.Item(99, 10).Value = 466.66666
$ws_ALC.Cells.Item(99, 12).Value = 1399.99998
$ws_ALC.Cells.Item(99, 14).Value = -4395.999980000001
$ws_ALC.Cells.Item(100, 8).Value = 2444.8
$ws_ALC.Cells.Item(100, 9).Value = 2556
$ws_ALC.Cells.Item(100, 10).Value = 2000
$ws_ALC.Cells.Item(100, 11).Value = 2556
$ws_ALC.Cells.Item(100, 12).Value = 2000
$ws_ALC.Cells.Item(100, 13).Value = -2015
$ws_ALC.Cells.Item(100, 14).Value = -3082
$ws_ALC.Cells.Item(102, 8).Value = 13428.286
$ws_ALC.Cells.Item(102, 10).Value = 13428.286
$ws_ALC.Cells.Item(102, 12).Value = 13428.286
$ws_ALC.Cells.Item(102, 14).Value = -19918.286
$ws_ALC.Cells.Item(105, 8).Value = 32921.715
$ws_ALC.Cells.Item(105, 10).Value = 32921.715
$ws_ALC.Cells.Item(105, 12).Value = 32921.715
$ws_ALC.Cells.Item(105, 14).Value = -39909.715
$ws_ALC.Cells.Item(106, 8).Value = 55558704
$ws_ALC.Cells.Item(106, 9).Value = 58826424
$ws_ALC.Cells.Item(106, 11).Value = 58826424
$ws_ALC.Cells.Item(106, 13).Value = -58825793
$ws_ALC.Cells.Item(107, 8).Value = 810.9231
$ws_ALC.Cells.Item(107, 9).Value = 759.2
$ws_ALC.Cells.Item(107, 11).Value = 759.2
$ws_ALC.Cells.Item(107, 13).Value = 1160.8
$ws_ALC.Cells.Item(111, 8).Value = 3418.75
$ws_ALC.Cells.Item(111, 9).Value = 3418.75
$ws_ALC.Cells.Item(111, 11).Value = 10256.25
$ws_ALC.Cells.Item(111, 13).Value = -7189.25
$ws_ALC.Cells.Item(112, 8).Value = 4153.88
$ws_ALC.Cells.Item(112, 10).Value = 4370.3184
$ws_ALC.Cells.Item(112, 12).Value = 13110.9552
$ws_ALC.Cells.Item(112, 14).Value = -15326.9552
$ws_ALC.Cells.Item(113, 8).Value = 7597
$ws_ALC.Cells.Item(113, 10).Value = 9328.333000000001
$ws_ALC.Cells.Item(113, 12).Value = 9328.333000000001
$ws_ALC.Cells.Item(113, 14).Value = -15836.333
$ws_ALC.Cells.Item(127, 8).Value = 1316.5714
$ws_ALC.Cells.Item(127, 9).Value = 1316.5714
$ws_ALC.Cells.Item(127, 11).Value = 3949.7142
$ws_ALC.Cells.Item(127, 13).Value = 1010.2858
$ws_ALC.Cells.Item(129, 8).Value = 1330.85
$ws_ALC.Cells.Item(129, 9).Value = 1038.5625
$ws_ALC.Cells.Item(129, 11).Value = 3115.6875
$ws_ALC.Cells.Item(129, 13).Value = 1884.3125
$ws_ALC.Cells.Item(132, 8).Value = 25519.453
$ws_ALC.Cells.Item(132, 9).Value = 31685.424
$ws_ALC.Cells.Item(132, 10).Value = 2910.889
$ws_ALC.Cells.Item(132, 11).Value = 95056.272
$ws_ALC.Cells.Item(132, 12).Value = 8732.667000000001
$ws_ALC.Cells.Item(132, 13).Value = -92526.272
$ws_ALC.Cells.Item(132, 14).Value = -13792.667
$ws_ALC.Cells.Item(135, 8).Value = 2941.1667
$ws_ALC.Cells.Item(135, 9).Value = 2945.8
$ws_ALC.Cells.Item(135, 10).Value = 2937.8572
$ws_ALC.Cells.Item(135, 11).Value = 26512.2
$ws_ALC.Cells.Item(135, 12).Value = 26440.7148
$ws_ALC.Cells.Item(135, 13).Value = -23977.2
$ws_ALC.Cells.Item(135, 14).Value = -31510.7148
$ws_ALC.Cells.Item(138, 8).Value = 6119.9585
$ws_ALC.Cells.Item(138, 9).Value = 13229
$ws_ALC.Cells.Item(138, 10).Value = 4249.1577
$ws_ALC.Cells.Item(138, 11).Value = 39687
$ws_ALC.Cells.Item(138, 12).Value = 12747.4731
$ws_ALC.Cells.Item(138, 13).Value = -34547
$ws_ALC.Cells.Item(138, 14).Value = -23027.4731
$ws_ALC.Cells.Item(141, 8).Value = 2536.9092
$ws_ALC.Cells.Item(141, 9).Value = 2536.9092
$ws_ALC.Cells.Item(141, 11).Value = 7610.7276
$ws_ALC.Cells.Item(141, 13).Value = -2430.7276
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(2, 8).Value = 1279.75
$ws_ARM.Cells.Item(2, 9).Value = 831.5294
$ws_ARM.Cells.Item(2, 10).Value = 3819.6667
$ws_ARM.Cells.Item(2, 11).Value = 831.5294
$ws_ARM.Cells.Item(2, 12).Value = 3819.6667
$ws_ARM.Cells.Item(2, 13).Value = -718.5294
$ws_ARM.Cells.Item(2, 14).Value = -4045.6667
$ws_ARM.Cells.Item(32, 8).Value = 229770.14
$ws_ARM.Cells.Item(32, 9).Value = 272110.12
$ws_ARM.Cells.Item(32, 10).Value = 5973
$ws_ARM.Cells.Item(32, 11).Value = 272110.12
$ws_ARM.Cells.Item(32, 12).Value = 5973
$ws_ARM.Cells.Item(32, 13).Value = -271823.12
$ws_ARM.Cells.Item(32, 14).Value = -6547
$ws_ARM.Cells.Item(61, 8).Value = 1371561.5
$ws_ARM.Cells.Item(61, 9).Value = 3527.647
$ws_ARM.Cells.Item(61, 10).Value = 6539689.5
$ws_ARM.Cells.Item(61, 11).Value = 3527.647
$ws_ARM.Cells.Item(61, 12).Value = 6539689.5
$ws_ARM.Cells.Item(61, 13).Value = -3315.647
$ws_ARM.Cells.Item(61, 14).Value = -6540113.5
$ws_ARM.Cells.Item(106, 8).Value = 0
$ws_ARM.Cells.Item(106, 10).Value = 0
$ws_ARM.Cells.Item(106, 12).Value = 0
$ws_ARM.Cells.Item(106, 14).ClearContents()
$ws_ARM.Cells.Item(113, 8).Value = 0
$ws_ARM.Cells.Item(113, 10).Value = 0
$ws_ARM.Cells.Item(113, 12).Value = 0
$ws_ARM.Cells.Item(113, 14).ClearContents()
$ws_ARM.Cells.Item(115, 8).Value = 0
$ws_ARM.Cells.Item(115, 10).Value = 0
$ws_ARM.Cells.Item(115, 12).Value = 0
$ws_ARM.Cells.Item(115, 14).ClearContents()
$ws_ARM.Cells.Item(116, 8).Value = 1279.75
$ws_ARM.Cells.Item(116, 9).Value = 831.5294
$ws_ARM.Cells.Item(116, 10).Value = 3819.6667
$ws_ARM.Cells.Item(116, 11).Value = 831.5294
$ws_ARM.Cells.Item(116, 12).Value = 3819.6667
$ws_ARM.Cells.Item(116, 13).Value = 1462.4706
$ws_ARM.Cells.Item(116, 14).Value = -8407.6667
$ws_ARM.Cells.Item(122, 8).Value = 1665.6578
$ws_ARM.Cells.Item(122, 9).Value = 1552.32
$ws_ARM.Cells.Item(122, 10).Value = 1883.6154
$ws_ARM.Cells.Item(122, 11).Value = 4656.96
$ws_ARM.Cells.Item(122, 12).Value = 5650.8462
$ws_ARM.Cells.Item(122, 13).Value = -2206.96
$ws_ARM.Cells.Item(122, 14).Value = -10550.8462
$ws_ARM.Cells.Item(131, 8).Value = 0
$ws_ARM.Cells.Item(131, 10).Value = 0
$ws_ARM.Cells.Item(131, 12).Value = 0
$ws_ARM.Cells.Item(131, 14).ClearContents()
$ws_ARM.Cells.Item(132, 8).Value = 1959.8837
$ws_ARM.Cells.Item(132, 9).Value = 1496.8182
$ws_ARM.Cells.Item(132, 10).Value = 3488
$ws_ARM.Cells.Item(132, 11).Value = 4490.4546
$ws_ARM.Cells.Item(132, 12).Value = 10464
$ws_ARM.Cells.Item(132, 13).Value = -1960.4546
$ws_ARM.Cells.Item(132, 14).Value = -15524
$ws_ARM.Cells.Item(136, 8).Value = 1371561.5
$ws_ARM.Cells.Item(136, 9).Value = 3527.647
$ws_ARM.Cells.Item(136, 10).Value = 6539689.5
$ws_ARM.Cells.Item(136, 11).Value = 10582.941
$ws_ARM.Cells.Item(136, 12).Value = 19619068.5
$ws_ARM.Cells.Item(136, 13).Value = -8032.940999999999
$ws_ARM.Cells.Item(136, 14).Value = -19624168.5
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(3, 8).Value = 1279.75
$ws_BSM.Cells.Item(3, 9).Value = 831.5294
$ws_BSM.Cells.Item(3, 10).Value = 3819.6667
$ws_BSM.Cells.Item(3, 11).Value = 831.5294
$ws_BSM.Cells.Item(3, 12).Value = 3819.6667
$ws_BSM.Cells.Item(3, 13).Value = -717.5294
$ws_BSM.Cells.Item(3, 14).Value = -4047.6667
$ws_BSM.Cells.Item(86, 8).Value = 4498.375
$ws_BSM.Cells.Item(86, 9).Value = 1497.5
$ws_BSM.Cells.Item(86, 10).Value = 7499.25
$ws_BSM.Cells.Item(86, 11).Value = 1497.5
$ws_BSM.Cells.Item(86, 12).Value = 7499.25
$ws_BSM.Cells.Item(86, 13).Value = -374.5
$ws_BSM.Cells.Item(86, 14).Value = -9745.25
$ws_BSM.Cells.Item(89, 8).Value = 4498.375
$ws_BSM.Cells.Item(89, 9).Value = 1497.5
$ws_BSM.Cells.Item(89, 10).Value = 7499.25
$ws_BSM.Cells.Item(89, 11).Value = 7487.5
$ws_BSM.Cells.Item(89, 12).Value = 37496.25
$ws_BSM.Cells.Item(89, 13).Value = -1871.5
$ws_BSM.Cells.Item(89, 14).Value = -48728.25
$ws_BSM.Cells.Item(105, 8).Value = 6060.9614
$ws_BSM.Cells.Item(105, 9).Value = 7643.1875
$ws_BSM.Cells.Item(105, 10).Value = 3529.4
$ws_BSM.Cells.Item(105, 11).Value = 7643.1875
$ws_BSM.Cells.Item(105, 12).Value = 3529.4
$ws_BSM.Cells.Item(105, 13).Value = -5896.1875
$ws_BSM.Cells.Item(105, 14).Value = -7023.4
$ws_BSM.Cells.Item(107, 8).Value = 5768.75
$ws_BSM.Cells.Item(107, 9).Value = 6471.968
$ws_BSM.Cells.Item(107, 10).Value = 3346.5557
$ws_BSM.Cells.Item(107, 11).Value = 6471.968
$ws_BSM.Cells.Item(107, 12).Value = 3346.5557
$ws_BSM.Cells.Item(107, 13).Value = -4551.968
$ws_BSM.Cells.Item(107, 14).Value = -7186.5557
$ws_BSM.Cells.Item(134, 8).Value = 27274866
$ws_BSM.Cells.Item(134, 9).Value = 1872.6666
$ws_BSM.Cells.Item(134, 10).Value = 100002850
$ws_BSM.Cells.Item(134, 11).Value = 5617.9998
$ws_BSM.Cells.Item(134, 12).Value = 300008550
$ws_BSM.Cells.Item(134, 13).Value = -3082.9998
$ws_BSM.Cells.Item(134, 14).Value = -300013620
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(58, 8).Value = 2688.3914
$ws_CRP.Cells.Item(58, 9).Value = 2300.5625
$ws_CRP.Cells.Item(58, 10).Value = 3574.8572
$ws_CRP.Cells.Item(58, 11).Value = 2300.5625
$ws_CRP.Cells.Item(58, 12).Value = 3574.8572
$ws_CRP.Cells.Item(58, 13).Value = -2097.5625
$ws_CRP.Cells.Item(58, 14).Value = -3980.8572
$ws_CRP.Cells.Item(62, 8).Value = 4593.923
$ws_CRP.Cells.Item(62, 9).Value = 4393.4165
$ws_CRP.Cells.Item(62, 11).Value = 4393.4165
$ws_CRP.Cells.Item(62, 13).Value = -3769.4165
$ws_CRP.Cells.Item(65, 8).Value = 4593.923
$ws_CRP.Cells.Item(65, 9).Value = 4393.4165
$ws_CRP.Cells.Item(65, 11).Value = 21967.0825
$ws_CRP.Cells.Item(65, 13).Value = -18847.0825
$ws_CRP.Cells.Item(86, 8).Value = 12280.9
$ws_CRP.Cells.Item(86, 9).Value = 7648.4
$ws_CRP.Cells.Item(86, 10).Value = 26178.4
$ws_CRP.Cells.Item(86, 11).Value = 7648.4
$ws_CRP.Cells.Item(86, 12).Value = 26178.4
$ws_CRP.Cells.Item(86, 13).Value = -6525.4
$ws_CRP.Cells.Item(86, 14).Value = -28424.4
$ws_CRP.Cells.Item(89, 8).Value = 12280.9
$ws_CRP.Cells.Item(89, 9).Value = 7648.4
$ws_CRP.Cells.Item(89, 10).Value = 26178.4
$ws_CRP.Cells.Item(89, 11).Value = 38242
$ws_CRP.Cells.Item(89, 12).Value = 130892
$ws_CRP.Cells.Item(89, 13).Value = -32626
$ws_CRP.Cells.Item(89, 14).Value = -142124
$ws_CRP.Cells.Item(92, 8).Value = 44999
$ws_CRP.Cells.Item(92, 10).Value = 44999
$ws_CRP.Cells.Item(92, 12).Value = 44999
$ws_CRP.Cells.Item(92, 14).Value = -49991
$ws_CRP.Cells.Item(94, 8).Value = 1968.3846
$ws_CRP.Cells.Item(94, 10).Value = 2054.4443
$ws_CRP.Cells.Item(94, 12).Value = 2054.4443
$ws_CRP.Cells.Item(94, 14).Value = -2956.4443
$ws_CRP.Cells.Item(100, 8).Value = 21088.092
$ws_CRP.Cells.Item(100, 10).Value = 21088.092
$ws_CRP.Cells.Item(100, 12).Value = 21088.092
$ws_CRP.Cells.Item(100, 14).Value = -23252.092
$ws_CRP.Cells.Item(107, 8).Value = 1255.5416
$ws_CRP.Cells.Item(107, 9).Value = 1073.8235
$ws_CRP.Cells.Item(107, 11).Value = 1073.8235
$ws_CRP.Cells.Item(107, 13).Value = 846.1765
$ws_CRP.Cells.Item(112, 8).Value = 68871.75
$ws_CRP.Cells.Item(112, 10).Value = 68871.75
$ws_CRP.Cells.Item(112, 12).Value = 68871.75
$ws_CRP.Cells.Item(112, 14).Value = -71825.75
$ws_CRP.Cells.Item(132, 8).Value = 2518.4583
$ws_CRP.Cells.Item(132, 9).Value = 2129.2632
$ws_CRP.Cells.Item(132, 10).Value = 3997.4
$ws_CRP.Cells.Item(132, 11).Value = 6387.7896
$ws_CRP.Cells.Item(132, 12).Value = 11992.2
$ws_CRP.Cells.Item(132, 13).Value = -3857.7896
$ws_CRP.Cells.Item(132, 14).Value = -17052.2
$ws_CRP.Cells.Item(136, 8).Value = 2688.3914
$ws_CRP.Cells.Item(136, 9).Value = 2300.5625
$ws_CRP.Cells.Item(136, 10).Value = 3574.8572
$ws_CRP.Cells.Item(136, 11).Value = 6901.6875
$ws_CRP.Cells.Item(136, 12).Value = 10724.5716
$ws_CRP.Cells.Item(136, 13).Value = -4351.6875
$ws_CRP.Cells.Item(136, 14).Value = -15824.5716
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(2, 8).Value = 385.94736
$ws_CUL.Cells.Item(2, 10).Value = 829.6
$ws_CUL.Cells.Item(2, 12).Value = 4977.6
$ws_CUL.Cells.Item(2, 14).Value = -5203.6
$ws_CUL.Cells.Item(7, 8).Value = 0
$ws_CUL.Cells.Item(7, 9).Value = 0
$ws_CUL.Cells.Item(7, 10).Value = 0
$ws_CUL.Cells.Item(7, 11).Value = 0
$ws_CUL.Cells.Item(7, 12).Value = 0
$ws_CUL.Cells.Item(7, 13).ClearContents()
$ws_CUL.Cells.Item(7, 14).ClearContents()
$ws_CUL.Cells.Item(33, 8).Value = 59564.883
$ws_CUL.Cells.Item(33, 9).Value = 244.33333
$ws_CUL.Cells.Item(33, 10).Value = 126300.5
$ws_CUL.Cells.Item(33, 11).Value = 1465.99998
$ws_CUL.Cells.Item(33, 12).Value = 757803
$ws_CUL.Cells.Item(33, 13).Value = -1182.99998
$ws_CUL.Cells.Item(33, 14).Value = -758369
$ws_CUL.Cells.Item(37, 8).Value = 77762
$ws_CUL.Cells.Item(37, 10).Value = 77762
$ws_CUL.Cells.Item(37, 12).Value = 233286
$ws_CUL.Cells.Item(37, 14).Value = -233510
$ws_CUL.Cells.Item(98, 8).Value = 357.73334
$ws_CUL.Cells.Item(98, 10).Value = 357.73334
$ws_CUL.Cells.Item(98, 12).Value = 1073.20002
$ws_CUL.Cells.Item(98, 14).Value = -4069.20002
$ws_CUL.Cells.Item(118, 8).Value = 8000
$ws_CUL.Cells.Item(118, 9).Value = 8000
$ws_CUL.Cells.Item(118, 11).Value = 24000
$ws_CUL.Cells.Item(118, 13).Value = -22757
$ws_CUL.Cells.Item(132, 8).Value = 1407.8462
$ws_CUL.Cells.Item(132, 9).Value = 1144.2222
$ws_CUL.Cells.Item(132, 11).Value = 10297.9998
$ws_CUL.Cells.Item(132, 13).Value = -7767.9998
$ws_CUL.Cells.Item(136, 8).Value = 5328.1665
$ws_CUL.Cells.Item(136, 9).Value = 693
$ws_CUL.Cells.Item(136, 10).Value = 9963.333000000001
$ws_CUL.Cells.Item(136, 11).Value = 2079
$ws_CUL.Cells.Item(136, 12).Value = 29889.999
$ws_CUL.Cells.Item(136, 13).Value = 3021
$ws_CUL.Cells.Item(136, 14).Value = -40089.999
$ws_CUL.Cells.Item(139, 8).Value = 3908209.2
$ws_CUL.Cells.Item(139, 9).Value = 5001068
$ws_CUL.Cells.Item(139, 11).Value = 15003204
$ws_CUL.Cells.Item(139, 13).Value = -14998064
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(80, 8).Value = 12639200
$ws_GSM.Cells.Item(80, 10).Value = 41810590
$ws_GSM.Cells.Item(80, 12).Value = 41810590
$ws_GSM.Cells.Item(80, 14).Value = -41812586
$ws_GSM.Cells.Item(83, 8).Value = 12639200
$ws_GSM.Cells.Item(83, 10).Value = 41810590
$ws_GSM.Cells.Item(83, 12).Value = 209052950
$ws_GSM.Cells.Item(83, 14).Value = -209062934
$ws_GSM.Cells.Item(97, 8).Value = 5142.2173
$ws_GSM.Cells.Item(97, 9).Value = 847.5
$ws_GSM.Cells.Item(97, 10).Value = 14958.714
$ws_GSM.Cells.Item(97, 11).Value = 847.5
$ws_GSM.Cells.Item(97, 12).Value = 14958.714
$ws_GSM.Cells.Item(97, 13).Value = -351.5
$ws_GSM.Cells.Item(97, 14).Value = -15950.714
$ws_GSM.Cells.Item(99, 8).Value = 11321.1
$ws_GSM.Cells.Item(99, 9).Value = 8690.223
$ws_GSM.Cells.Item(99, 11).Value = 8690.223
$ws_GSM.Cells.Item(99, 13).Value = -6444.223
$ws_GSM.Cells.Item(132, 8).Value = 628266.7
$ws_GSM.Cells.Item(132, 9).Value = 2301.24
$ws_GSM.Cells.Item(132, 10).Value = 1606337.8
$ws_GSM.Cells.Item(132, 11).Value = 6903.719999999999
$ws_GSM.Cells.Item(132, 12).Value = 4819013.4
$ws_GSM.Cells.Item(132, 13).Value = -4373.719999999999
$ws_GSM.Cells.Item(132, 14).Value = -4824073.4
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(3, 8).Value = 10000
$ws_LTW.Cells.Item(3, 10).Value = 10000
$ws_LTW.Cells.Item(3, 12).Value = 10000
$ws_LTW.Cells.Item(3, 14).Value = -10224
$ws_LTW.Cells.Item(4, 8).Value = 19899
$ws_LTW.Cells.Item(4, 10).Value = 19899
$ws_LTW.Cells.Item(4, 12).Value = 19899
$ws_LTW.Cells.Item(4, 14).Value = -20125
$ws_LTW.Cells.Item(5, 8).Value = 31999.5
$ws_LTW.Cells.Item(5, 10).Value = 31999.5
$ws_LTW.Cells.Item(5, 12).Value = 31999.5
$ws_LTW.Cells.Item(5, 14).Value = -32225.5
$ws_LTW.Cells.Item(7, 8).Value = 7889.7827
$ws_LTW.Cells.Item(7, 9).Value = 4133.6
$ws_LTW.Cells.Item(7, 10).Value = 14932.625
$ws_LTW.Cells.Item(7, 11).Value = 4133.6
$ws_LTW.Cells.Item(7, 12).Value = 14932.625
$ws_LTW.Cells.Item(7, 13).Value = -4021.6
$ws_LTW.Cells.Item(7, 14).Value = -15156.625
$ws_LTW.Cells.Item(15, 8).Value = 10000
$ws_LTW.Cells.Item(15, 10).Value = 10000
$ws_LTW.Cells.Item(15, 12).Value = 10000
$ws_LTW.Cells.Item(15, 14).Value = -10340
$ws_LTW.Cells.Item(16, 8).Value = 3136.12
$ws_LTW.Cells.Item(16, 9).Value = 2876.3333
$ws_LTW.Cells.Item(16, 11).Value = 2876.3333
$ws_LTW.Cells.Item(16, 13).Value = -2706.3333
$ws_LTW.Cells.Item(28, 8).Value = 19899
$ws_LTW.Cells.Item(28, 10).Value = 19899
$ws_LTW.Cells.Item(28, 12).Value = 19899
$ws_LTW.Cells.Item(28, 14).Value = -20363
$ws_LTW.Cells.Item(37, 8).Value = 19899
$ws_LTW.Cells.Item(37, 10).Value = 19899
$ws_LTW.Cells.Item(37, 12).Value = 19899
$ws_LTW.Cells.Item(37, 14).Value = -20113
$ws_LTW.Cells.Item(60, 8).Value = 0
$ws_LTW.Cells.Item(60, 10).Value = 0
$ws_LTW.Cells.Item(60, 12).Value = 0
$ws_LTW.Cells.Item(60, 14).ClearContents()
$ws_LTW.Cells.Item(61, 8).Value = 4361.0557
$ws_LTW.Cells.Item(61, 9).Value = 5241.4
$ws_LTW.Cells.Item(61, 10).Value = 3260.625
$ws_LTW.Cells.Item(61, 11).Value = 5241.4
$ws_LTW.Cells.Item(61, 12).Value = 3260.625
$ws_LTW.Cells.Item(61, 13).Value = -5039.4
$ws_LTW.Cells.Item(61, 14).Value = -3664.625
$ws_LTW.Cells.Item(74, 8).Value = 41304.25
$ws_LTW.Cells.Item(74, 9).Value = 32611
$ws_LTW.Cells.Item(74, 10).Value = 49997.5
$ws_LTW.Cells.Item(74, 11).Value = 32611
$ws_LTW.Cells.Item(74, 12).Value = 49997.5
$ws_LTW.Cells.Item(74, 13).Value = -31613
$ws_LTW.Cells.Item(74, 14).Value = -51993.5
$ws_LTW.Cells.Item(77, 8).Value = 41304.25
$ws_LTW.Cells.Item(77, 9).Value = 32611
$ws_LTW.Cells.Item(77, 10).Value = 49997.5
$ws_LTW.Cells.Item(77, 11).Value = 97833
$ws_LTW.Cells.Item(77, 12).Value = 149992.5
$ws_LTW.Cells.Item(77, 13).Value = -92841
$ws_LTW.Cells.Item(77, 14).Value = -159976.5
$ws_LTW.Cells.Item(95, 8).Value = 36172
$ws_LTW.Cells.Item(95, 10).Value = 36172
$ws_LTW.Cells.Item(95, 12).Value = 36172
$ws_LTW.Cells.Item(95, 14).Value = -41664
$ws_LTW.Cells.Item(100, 8).Value = 3345.1428
$ws_LTW.Cells.Item(100, 9).Value = 2488.4666
$ws_LTW.Cells.Item(100, 10).Value = 5486.8335
$ws_LTW.Cells.Item(100, 11).Value = 2488.4666
$ws_LTW.Cells.Item(100, 12).Value = 5486.8335
$ws_LTW.Cells.Item(100, 13).Value = -1947.4666
$ws_LTW.Cells.Item(100, 14).Value = -6568.8335
$ws_LTW.Cells.Item(113, 8).Value = 4361.0557
$ws_LTW.Cells.Item(113, 9).Value = 5241.4
$ws_LTW.Cells.Item(113, 10).Value = 3260.625
$ws_LTW.Cells.Item(113, 11).Value = 5241.4
$ws_LTW.Cells.Item(113, 12).Value = 3260.625
$ws_LTW.Cells.Item(113, 13).Value = -3071.4
$ws_LTW.Cells.Item(113, 14).Value = -7600.625
$ws_LTW.Cells.Item(122, 8).Value = 3114.3948
$ws_LTW.Cells.Item(122, 9).Value = 2526.5386
$ws_LTW.Cells.Item(122, 10).Value = 4388.0835
$ws_LTW.Cells.Item(122, 11).Value = 7579.6158
$ws_LTW.Cells.Item(122, 12).Value = 13164.2505
$ws_LTW.Cells.Item(122, 13).Value = -5129.6158
$ws_LTW.Cells.Item(122, 14).Value = -18064.2505
$ws_LTW.Cells.Item(126, 8).Value = 7889.7827
$ws_LTW.Cells.Item(126, 9).Value = 4133.6
$ws_LTW.Cells.Item(126, 10).Value = 14932.625
$ws_LTW.Cells.Item(126, 11).Value = 12400.8
$ws_LTW.Cells.Item(126, 12).Value = 44797.875
$ws_LTW.Cells.Item(126, 13).Value = -9930.800000000001
$ws_LTW.Cells.Item(126, 14).Value = -49737.875
$ws_LTW.Cells.Item(130, 8).Value = 177776.5
$ws_LTW.Cells.Item(130, 10).Value = 177776.5
$ws_LTW.Cells.Item(130, 12).Value = 177776.5
$ws_LTW.Cells.Item(130, 14).Value = -187816.5
$ws_LTW.Cells.Item(132, 8).Value = 2740.8948
$ws_LTW.Cells.Item(132, 9).Value = 2671.9443
$ws_LTW.Cells.Item(132, 11).Value = 8015.8329
$ws_LTW.Cells.Item(132, 13).Value = -5485.8329
$ws_LTW.Cells.Item(136, 8).Value = 3411.7017
$ws_LTW.Cells.Item(136, 9).Value = 2222.3704
$ws_LTW.Cells.Item(136, 10).Value = 4482.1
$ws_LTW.Cells.Item(136, 11).Value = 6667.111199999999
$ws_LTW.Cells.Item(136, 12).Value = 13446.3
$ws_LTW.Cells.Item(136, 13).Value = -4117.111199999999
$ws_LTW.Cells.Item(136, 14).Value = -18546.3
$ws_LTW.Cells.Item(137, 8).Value = 123326.664
$ws_LTW.Cells.Item(137, 10).Value = 123326.664
$ws_LTW.Cells.Item(137, 12).Value = 123326.664
$ws_LTW.Cells.Item(137, 14).Value = -133526.664
$ws_LTW.Cells.Item(140, 8).Value = 111249.25
$ws_LTW.Cells.Item(140, 10).Value = 111249.25
$ws_LTW.Cells.Item(140, 12).Value = 111249.25
$ws_LTW.Cells.Item(140, 14).Value = -121609.25
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(19, 8).Value = 102235.336
$ws_WVR.Cells.Item(19, 10).Value = 152753
$ws_WVR.Cells.Item(19, 12).Value = 152753
$ws_WVR.Cells.Item(19, 14).Value = -153101
$ws_WVR.Cells.Item(46, 8).Value = 51867.91
$ws_WVR.Cells.Item(46, 10).Value = 51867.91
$ws_WVR.Cells.Item(46, 12).Value = 51867.91
$ws_WVR.Cells.Item(46, 14).Value = -52329.91
$ws_WVR.Cells.Item(62, 8).Value = 4299.5
$ws_WVR.Cells.Item(62, 9).Value = 3899.3333
$ws_WVR.Cells.Item(62, 11).Value = 3899.3333
$ws_WVR.Cells.Item(62, 13).Value = -3275.3333
$ws_WVR.Cells.Item(65, 8).Value = 4299.5
$ws_WVR.Cells.Item(65, 9).Value = 3899.3333
$ws_WVR.Cells.Item(65, 11).Value = 19496.6665
$ws_WVR.Cells.Item(65, 13).Value = -16376.6665
$ws_WVR.Cells.Item(81, 8).Value = 2626.7334
$ws_WVR.Cells.Item(81, 9).Value = 3209.6
$ws_WVR.Cells.Item(81, 10).Value = 1461
$ws_WVR.Cells.Item(81, 11).Value = 6419.2
$ws_WVR.Cells.Item(81, 12).Value = 2922
$ws_WVR.Cells.Item(81, 13).Value = -5358.2
$ws_WVR.Cells.Item(81, 14).Value = -5044
$ws_WVR.Cells.Item(84, 8).Value = 2626.7334
$ws_WVR.Cells.Item(84, 9).Value = 3209.6
$ws_WVR.Cells.Item(84, 10).Value = 1461
$ws_WVR.Cells.Item(84, 11).Value = 32096
$ws_WVR.Cells.Item(84, 12).Value = 14610
$ws_WVR.Cells.Item(84, 13).Value = -26792
$ws_WVR.Cells.Item(84, 14).Value = -25218
$ws_WVR.Cells.Item(95, 8).Value = 40498.5
$ws_WVR.Cells.Item(95, 9).Value = 25000
$ws_WVR.Cells.Item(95, 10).Value = 45664.668
$ws_WVR.Cells.Item(95, 11).Value = 25000
$ws_WVR.Cells.Item(95, 12).Value = 45664.668
$ws_WVR.Cells.Item(95, 13).Value = -22254
$ws_WVR.Cells.Item(95, 14).Value = -51156.668
$ws_WVR.Cells.Item(107, 8).Value = 1058870.5
$ws_WVR.Cells.Item(107, 9).Value = 683.3182
$ws_WVR.Cells.Item(107, 10).Value = 5714894
$ws_WVR.Cells.Item(107, 11).Value = 2049.9546
$ws_WVR.Cells.Item(107, 12).Value = 17144682
$ws_WVR.Cells.Item(107, 13).Value = -129.9546
$ws_WVR.Cells.Item(107, 14).Value = -17148522
$ws_WVR.Cells.Item(112, 8).Value = 30791.5
$ws_WVR.Cells.Item(112, 10).Value = 30791.5
$ws_WVR.Cells.Item(112, 12).Value = 30791.5
$ws_WVR.Cells.Item(112, 14).Value = -33745.5
$ws_WVR.Cells.Item(122, 8).Value = 2148.2424
$ws_WVR.Cells.Item(122, 9).Value = 2096.4644
$ws_WVR.Cells.Item(122, 10).Value = 2438.2
$ws_WVR.Cells.Item(122, 11).Value = 6289.3932
$ws_WVR.Cells.Item(122, 12).Value = 7314.599999999999
$ws_WVR.Cells.Item(122, 13).Value = -3839.3932
$ws_WVR.Cells.Item(122, 14).Value = -12214.6
$ws_WVR.Cells.Item(125, 8).Value = 54881.645
$ws_WVR.Cells.Item(125, 10).Value = 54881.645
$ws_WVR.Cells.Item(125, 12).Value = 54881.645
$ws_WVR.Cells.Item(125, 14).Value = -64721.645
$ws_WVR.Cells.Item(130, 8).Value = 73499.5
$ws_WVR.Cells.Item(130, 10).Value = 73499.5
$ws_WVR.Cells.Item(130, 12).Value = 73499.5
$ws_WVR.Cells.Item(130, 14).Value = -83539.5
$ws_WVR.Cells.Item(132, 8).Value = 2771.3142
$ws_WVR.Cells.Item(132, 9).Value = 2600.08
$ws_WVR.Cells.Item(132, 11).Value = 7800.24
$ws_WVR.Cells.Item(132, 13).Value = -5270.24
$ws_WVR.Cells.Item(134, 8).Value = 51867.91
$ws_WVR.Cells.Item(134, 10).Value = 51867.91
$ws_WVR.Cells.Item(134, 12).Value = 155603.73
$ws_WVR.Cells.Item(134, 14).Value = -160673.73
$ws_WVR.Cells.Item(136, 8).Value = 32856.37
$ws_WVR.Cells.Item(136, 9).Value = 44800.8
$ws_WVR.Cells.Item(136, 10).Value = 2995.3
$ws_WVR.Cells.Item(136, 11).Value = 134402.4
$ws_WVR.Cells.Item(136, 12).Value = 8985.900000000001
$ws_WVR.Cells.Item(136, 13).Value = -131852.4
$ws_WVR.Cells.Item(136, 14).Value = -14085.9
